$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.651.69"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "2.580.00"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'581.22"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'145.00"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +1.44%  "

$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").Value = "'26.97"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").Value = "3.040.91"
$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("D15").Value = "62.547.34"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").Value = "2.576.61"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "'11.19"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").Value = "'338.23"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").Value = "'66.95"
$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("D24").Value = "2.699.98"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  -2.28%  "

$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.29%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'7.87"
$ws.Range("E28").Value = "  +3.22%  "

$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "'1.47"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "'461.11"
$ws.Range("E32").Value = "  +10.63%  "

$ws.Range("D33").Value = "0.0₃0808"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").Value = "'176.81"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "'18.89"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D41").Value = "'1.69"
$ws.Range("E41").Value = "  -3.02%  "

$ws.Range("D42").Value = "'157.22"
$ws.Range("E42").Value = "  +4.48%  "

$ws.Range("D43").Value = "'3.72"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").Value = "'21.13"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("D45").Value = "'0.627"
$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("D46").Value = "'0.0535"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").Value = "'18.08"
$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("E51").Value = "  -0.21%  "
